$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44553
$ws.Range("M2").Value = 200
$ws.Range("D3").Value = 44553
$ws.Range("M3").Value = 150
$ws.Range("D4").Value = 44558
$ws.Range("L4").Value = 'Especial'
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 22000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 22000
$ws.Range("S4").Value = 3667
$ws.Range("D5").Value = 44558
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = '$/bandeja 6 kilos'
$ws.Range("R5").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S5").Value = 3000
$ws.Range("T5").Value = 6
$ws.Range("D7").Value = 44187
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 14000
$ws.Range("O7").Value = 14000
$ws.Range("P7").Value = 14000
$ws.Range("Q7").Value = '$/bandeja 7 kilos'
$ws.Range("S7").Value = 2000
$ws.Range("T7").Value = 7
$ws.Range("D8").Value = 44187
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("Q8").Value = '$/bandeja 7 kilos'
$ws.Range("S8").Value = 1714
$ws.Range("T8").Value = 7
$ws.Range("D9").Value = 44561
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 18000
$ws.Range("O9").Value = 18000
$ws.Range("P9").Value = 18000
$ws.Range("R9").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S9").Value = 3000
$ws.Range("D12").Value = 44550
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 24000
$ws.Range("O12").Value = 24000
$ws.Range("P12").Value = 24000
$ws.Range("R12").Value = 'Región Metropolitana'
$ws.Range("S12").Value = 3429
$ws.Range("D13").Value = 44572
$ws.Range("M13").Value = 65
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("Q13").Value = '$/bandeja 6 kilos'
$ws.Range("R13").Value = 'Región Metropolitana'
$ws.Range("S13").Value = 3333
$ws.Range("T13").Value = 6
